$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '28.954.98'
$ws.Cells.Item(2, 5).Value = '  +5.51%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.913.29'
$ws.Cells.Item(3, 5).Value = '  +4.78%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.45%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '338.96'
$ws.Cells.Item(5, 5).Value = '  +2.16%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  -0.34%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4744'
$ws.Cells.Item(7, 5).Value = '  +3.71%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.4057'
$ws.Cells.Item(8, 5).Value = '  +6.82%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '48.06'
$ws.Cells.Item(9, 5).Value = '  +3.60%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.08170'
$ws.Cells.Item(10, 5).Value = '  +3.57%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.030'
$ws.Cells.Item(11, 5).Value = '  +6.33%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '22.45'
$ws.Cells.Item(12, 5).Value = '  +6.82%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '6.092'
$ws.Cells.Item(13, 5).Value = '  +3.77%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '1.888.90'
$ws.Cells.Item(14, 5).Value = '  +3.08%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.390'
$ws.Cells.Item(15, 5).Value = '  +4.74%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '91.64'
$ws.Cells.Item(16, 5).Value = '  +2.39%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  -0.40%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +2.85%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06633'
$ws.Cells.Item(19, 5).Value = '  -0.01%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.87'
$ws.Cells.Item(20, 5).Value = '  +4.67%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.001'
$ws.Cells.Item(21, 5).Value = '  -0.35%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '28.984.68'
$ws.Cells.Item(22, 5).Value = '  +5.63%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.576'
$ws.Cells.Item(23, 5).Value = '  +4.70%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.22'
$ws.Cells.Item(24, 5).Value = '  +3.95%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -1.10%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.127.61'
$ws.Cells.Item(26, 5).Value = '  +3.94%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '160.63'
$ws.Cells.Item(27, 5).Value = '  +3.23%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '20.07'
$ws.Cells.Item(28, 5).Value = '  +3.68%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.183'
$ws.Cells.Item(29, 5).Value = '  +6.06%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '5.535'
$ws.Cells.Item(30, 5).Value = '  +4.93%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '121.12'
$ws.Cells.Item(31, 5).Value = '  +2.31%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.012'
$ws.Cells.Item(32, 5).Value = '  +7.69%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.09596'
$ws.Cells.Item(33, 5).Value = '  +3.26%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'ARBITRUM'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.431'
$ws.Cells.Item(34, 5).Value = '  +7.83%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'HuobiToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.646'
$ws.Cells.Item(35, 5).Value = '  +1.70%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.429'
$ws.Cells.Item(36, 5).Value = '  +3.57%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.06208'
$ws.Cells.Item(37, 5).Value = '  +4.70%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.02284'
$ws.Cells.Item(38, 5).Value = '  +4.93%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '8.668'
$ws.Cells.Item(39, 5).Value = '  +7.59%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.208'
$ws.Cells.Item(40, 5).Value = '  +5.81%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.6044'
$ws.Cells.Item(41, 5).Value = '  +4.95%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Aptos'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '10.59'
$ws.Cells.Item(42, 5).Value = '  +6.29%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Algorand'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.1902'
$ws.Cells.Item(43, 5).Value = '  +4.24%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -0.33%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.284'
$ws.Cells.Item(45, 5).Value = '  +1.47%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '12.36'
$ws.Cells.Item(46, 5).Value = '  +3.88%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.5643'
$ws.Cells.Item(47, 5).Value = '  +3.69%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.982'
$ws.Cells.Item(48, 5).Value = '  +6.32%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.07268'
$ws.Cells.Item(49, 5).Value = '  +9.91%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.153'
$ws.Cells.Item(50, 5).Value = '  +19.84%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '113.45'
$ws.Cells.Item(51, 5).Value = '  +2.26%  '
